$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (강아지족) - remove commas from description
$ws.Range("G2").Value = "강아지는 그들이 검을 받았을 때 그 검의 명예를 찾아서 떠난다. 식탐이 많아서 식도락 여행으로 오인받기도 하지만... 일단은 숭고한 것이다."

# Row 3 (고양이족) - rich text: "보이지만, 실상은" -> "보이지만... 그냥" (with font change on "그냥")
$cell3 = $ws.Range("G3")
$cell3.Value = "고양이들은 빨라야 한다고 교육받는다. 그럼에도 그들의 천성은 느긋하다. 숭고한 여행으로 보이지만... 그냥 느긋하게 즐기고 있는 것이다."
$cell3.Characters(1, 5).Font.Name = "맑은 고딕"
$cell3.Characters(6, 1).Font.Name = "Arial"
$cell3.Characters(7, 3).Font.Name = "맑은 고딕"
$cell3.Characters(10, 1).Font.Name = "Arial"
$cell3.Characters(11, 3).Font.Name = "맑은 고딕"
$cell3.Characters(14, 1).Font.Name = "Arial"
$cell3.Characters(15, 5).Font.Name = "맑은 고딕"
$cell3.Characters(20, 2).Font.Name = "Arial"
$cell3.Characters(22, 4).Font.Name = "맑은 고딕"
$cell3.Characters(26, 1).Font.Name = "Arial"
$cell3.Characters(27, 3).Font.Name = "맑은 고딕"
$cell3.Characters(30, 1).Font.Name = "Arial"
$cell3.Characters(31, 3).Font.Name = "맑은 고딕"
$cell3.Characters(34, 1).Font.Name = "Arial"
$cell3.Characters(35, 4).Font.Name = "맑은 고딕"
$cell3.Characters(39, 2).Font.Name = "Arial"
$cell3.Characters(41, 3).Font.Name = "맑은 고딕"
$cell3.Characters(44, 1).Font.Name = "Arial"
$cell3.Characters(45, 4).Font.Name = "맑은 고딕"
$cell3.Characters(49, 1).Font.Name = "Arial"
$cell3.Characters(50, 7).Font.Name = "맑은 고딕"
$cell3.Characters(57, 1).Font.Name = "Arial"
$cell3.Characters(58, 2).Font.Name = "Arial Unicode MS"
$cell3.Characters(60, 1).Font.Name = "Arial"
$cell3.Characters(61, 4).Font.Name = "맑은 고딕"
$cell3.Characters(65, 1).Font.Name = "Arial"
$cell3.Characters(66, 3).Font.Name = "맑은 고딕"
$cell3.Characters(69, 1).Font.Name = "Arial"
$cell3.Characters(70, 6).Font.Name = "맑은 고딕"
$cell3.Characters(76, 1).Font.Name = "Arial"

# Row 5 (쥐돌이족) - description rewritten
$ws.Range("G5").Value = "그들은 미지를 동경하고 신비를 탐험하는 자들. 그들은 발견을 숭배하고 탐험을 지향한다. 음식에도 그러하다."

# Row 6 (수면양족) - remove comma
$ws.Range("G6").Value = "푸근해 보이는 외형과는 달리 한번 불타면 재가 될 때 까지 불태우는 종족. 당신은 모험에 불이 붙었다. 탐험이라는 연못으로 뛰어들어라."

# Update selection to match target
$ws.Range("G7").Select()
